$wb = $excel.ActiveWorkbook

# The edit replaces the Team ID and the four "Executed By" names in the
# "Shopenzer Testcases" sheet with a new set of values.
$ws = $wb.Worksheets.Item("Shopenzer Testcases")
$ws.Activate()

# Team ID value (row 2, column F)
$ws.Range("F2").Value = "PNT2022TMID53380"

# Executed By column (N6:N9)
$ws.Range("N6").Value = "Ritunjay M"
$ws.Range("N7").Value = "Praveen Raagul R"
$ws.Range("N8").Value = "Pradeep V"
$ws.Range("N9").Value = "Munish Kumar S"

# Reflect the final selection left behind in the saved file.
$ws.Range("N11").Select()
